# Applies the "Automatic update of files" change:
#   1) Column C ("Förändrad") on every data row (2-33) is bumped from the
#      old serial date to 45186 (2023-09-17).
#   2) On the rows that still carry the per-case document links (rows 2-5,
#      columns S/T/V/W/X/Y), the HYPERLINK() formulas gain a second
#      ("friendly name") argument equal to the case id in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Bump the "Förändrad" date (column C) for every data row -----------
for ($r = 2; $r -le 33; $r++) {
    $ws.Range("C$r").Value = 45186
}

# --- 2) Rewrite the HYPERLINK formulas to include the friendly name -------
$linkFolders = @{
    "S" = "artfynd"
    "T" = "kartor"
    "V" = "klagomål"
    "W" = "klagomålsmail"
    "X" = "tillsyn"
    "Y" = "tillsynsmail"
}
$linkExt = @{
    "S" = "xlsx"
    "T" = "png"
    "V" = "docx"
    "W" = "docx"
    "X" = "docx"
    "Y" = "docx"
}

for ($r = 2; $r -le 5; $r++) {
    $beteckning = $ws.Range("A$r").Value2
    foreach ($col in @("S", "T", "V", "W", "X", "Y")) {
        $folder = $linkFolders[$col]
        $ext = $linkExt[$col]
        $url = "https://klasma.github.io/Logging_MOLNDAL/$folder/$beteckning.$ext"
        $formula = '=HYPERLINK("' + $url + '", "' + $beteckning + '")'
        $ws.Range("$col$r").Formula = $formula
    }
}
